{"js": "// Remove the red \"(This is a change \u2013 Version for branch alternate)\" annotation\n// (and the two leading spaces before it) from the first paragraph, leaving only\n// \"This is a Microsoft word document.\", then drop the trailing empty paragraph\n// that sits right before the section break at the end of the body.\n\nconst body = context.document.body;\n\n// 1) Strip the \"  (This is a change \u2013 Version for branch alternate)\" text that\n//    was appended after \"This is a Microsoft word document.\"\nconst target = body.search(\"  (This is a change \\u2013 Version for branch alternate)\", { matchCase: true });\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  target.items[0].delete();\n  await context.sync();\n}\n\n// 2) Remove the trailing empty paragraph right before the final section break.\n//    A Range.delete() on an empty range at the very end of the body is a\n//    no-op (there is nothing \"selected\" to remove), so instead select from\n//    the start of the last-but-one paragraph through the end of the very\n//    last paragraph and delete that combined range \u2014 this removes the\n//    last-but-one paragraph's mark, merging the (empty) last paragraph away\n//    while the last-but-one paragraph's own properties/style survive.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst last = items[items.length - 1];\nif (last && last.text === \"\" && items.length >= 2) {\n  const secondLast = items[items.length - 2];\n  const combined = secondLast.getRange().expandTo(last.getRange());\n  combined.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the red \"(This is a change \u2013 Version for branch alternate)\" annotation\n# (and the two leading spaces before it) from the first paragraph, leaving only\n# \"This is a Microsoft word document.\", then drop the trailing empty paragraph\n# that sits right before the section break at the end of the body.\n\n$d = $word.ActiveDocument\n\n# 1) Strip the \"  (This is a change \u2013 Version for branch alternate)\" text that\n#    was appended after \"This is a Microsoft word document.\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Text = \"  (This is a change \u2013 Version for branch alternate)\"\n$rng.Find.Replacement.Text = \"\"\n$rng.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 2) Remove the trailing empty paragraph right before the final section break.\n#    A Range.Delete() on an empty range at the very end of the body is a\n#    no-op (there is nothing \"selected\" to remove), so instead select from\n#    the start of the last-but-one paragraph through the end of the very\n#    last paragraph and delete that combined range \u2014 this removes the\n#    last-but-one paragraph's mark, merging the (empty) last paragraph away\n#    while the last-but-one paragraph's own properties/style survive.\n$count = $d.Paragraphs.Count\n$last = $d.Paragraphs.Item($count)\nif ($last.Range.Text.TrimEnd([char]13, [char]7) -eq \"\" -and $count -ge 2) {\n    $secondLast = $d.Paragraphs.Item($count - 1)\n    $mergeRange = $d.Range($secondLast.Range.Start, $last.Range.End)\n    $mergeRange.Delete()\n}\n"}
